$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 2400  # ALC H6: 201.95 -> 2400
$ws.Cells.Item(6, 9).Value = 3766.6667  # ALC I6: 185.3077 -> 3766.6667
$ws.Cells.Item(6, 10).Value = 350  # ALC J6: 232.85715 -> 350
$ws.Cells.Item(6, 11).Value = 11300.0001  # ALC K6: 555.9231 -> 11300.0001
$ws.Cells.Item(6, 12).Value = 1050  # ALC L6: 698.5714499999999 -> 1050
$ws.Cells.Item(6, 13).Value = -11188.0001  # ALC M6: -443.9231 -> -11188.0001
$ws.Cells.Item(6, 14).Value = -1274  # ALC N6: -922.5714499999999 -> -1274

$ws.Cells.Item(76, 8).Value = 2962.5356  # ALC H76: 2890.3438 -> 2962.5356
$ws.Cells.Item(76, 9).Value = 2850.5264  # ALC I76: 2795.238 -> 2850.5264
$ws.Cells.Item(76, 10).Value = 3199  # ALC J76: 3071.9092 -> 3199
$ws.Cells.Item(76, 11).Value = 2850.5264  # ALC K76: 2795.238 -> 2850.5264
$ws.Cells.Item(76, 12).Value = 3199  # ALC L76: 3071.9092 -> 3199
$ws.Cells.Item(76, 13).Value = -2535.5264  # ALC M76: -2480.238 -> -2535.5264
$ws.Cells.Item(76, 14).Value = -3829  # ALC N76: -3701.9092 -> -3829

$ws.Cells.Item(79, 8).Value = 2962.5356  # ALC H79: 2890.3438 -> 2962.5356
$ws.Cells.Item(79, 9).Value = 2850.5264  # ALC I79: 2795.238 -> 2850.5264
$ws.Cells.Item(79, 10).Value = 3199  # ALC J79: 3071.9092 -> 3199
$ws.Cells.Item(79, 11).Value = 2850.5264  # ALC K79: 2795.238 -> 2850.5264
$ws.Cells.Item(79, 12).Value = 3199  # ALC L79: 3071.9092 -> 3199
$ws.Cells.Item(79, 13).Value = -1758.5264  # ALC M79: -1703.238 -> -1758.5264
$ws.Cells.Item(79, 14).Value = -5383  # ALC N79: -5255.9092 -> -5383

$ws.Cells.Item(127, 8).Value = 686.1  # ALC H127: 863.4286 -> 686.1
$ws.Cells.Item(127, 9).Value = 594.4286  # ALC I127: 654.8 -> 594.4286
$ws.Cells.Item(127, 10).Value = 900  # ALC J127: 979.3333 -> 900
$ws.Cells.Item(127, 11).Value = 1783.2858  # ALC K127: 1964.4 -> 1783.2858
$ws.Cells.Item(127, 12).Value = 2700  # ALC L127: 2937.9999 -> 2700
$ws.Cells.Item(127, 13).Value = 3176.7142  # ALC M127: 2995.6 -> 3176.7142
$ws.Cells.Item(127, 14).Value = -12620  # ALC N127: -12857.9999 -> -12620

$ws.Cells.Item(141, 8).Value = 3929.4048  # ALC H141: 4765.525 -> 3929.4048
$ws.Cells.Item(141, 9).Value = 1737.7368  # ALC I141: 2220 -> 1737.7368
$ws.Cells.Item(141, 10).Value = 24750.25  # ALC J141: 12402.1 -> 24750.25
$ws.Cells.Item(141, 11).Value = 5213.2104  # ALC K141: 6660 -> 5213.2104
$ws.Cells.Item(141, 12).Value = 74250.75  # ALC L141: 37206.3 -> 74250.75
$ws.Cells.Item(141, 13).Value = -33.21039999999994  # ALC M141: -1480 -> -33.21039999999994
$ws.Cells.Item(141, 14).Value = -84610.75  # ALC N141: -47566.3 -> -84610.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 12654.397  # ARM H32: 6923.8 -> 12654.397
$ws.Cells.Item(32, 9).Value = 3539.65  # ARM I32: 2846.3896 -> 3539.65
$ws.Cells.Item(32, 10).Value = 32186  # ARM J32: 20574.262 -> 32186
$ws.Cells.Item(32, 11).Value = 3539.65  # ARM K32: 2846.3896 -> 3539.65
$ws.Cells.Item(32, 12).Value = 32186  # ARM L32: 20574.262 -> 32186
$ws.Cells.Item(32, 13).Value = -3252.65  # ARM M32: -2559.3896 -> -3252.65
$ws.Cells.Item(32, 14).Value = -32760  # ARM N32: -21148.262 -> -32760

$ws.Cells.Item(61, 8).Value = 1152.303  # ARM H61: 1208.4667 -> 1152.303
$ws.Cells.Item(61, 9).Value = 1003.9048  # ARM I61: 1098.8235 -> 1003.9048
$ws.Cells.Item(61, 10).Value = 1412  # ARM J61: 1351.8462 -> 1412
$ws.Cells.Item(61, 11).Value = 1003.9048  # ARM K61: 1098.8235 -> 1003.9048
$ws.Cells.Item(61, 12).Value = 1412  # ARM L61: 1351.8462 -> 1412
$ws.Cells.Item(61, 13).Value = -791.9048  # ARM M61: -886.8235 -> -791.9048
$ws.Cells.Item(61, 14).Value = -1836  # ARM N61: -1775.8462 -> -1836

$ws.Cells.Item(110, 8).Value = 3798.0208  # ARM H110: 5109.5557 -> 3798.0208
$ws.Cells.Item(110, 9).Value = 4441.5586  # ARM I110: 6492.696 -> 4441.5586
$ws.Cells.Item(110, 10).Value = 2235.1428  # ARM J110: 2662.4614 -> 2235.1428
$ws.Cells.Item(110, 11).Value = 4441.5586  # ARM K110: 6492.696 -> 4441.5586
$ws.Cells.Item(110, 12).Value = 2235.1428  # ARM L110: 2662.4614 -> 2235.1428
$ws.Cells.Item(110, 13).Value = -2396.5586  # ARM M110: -4447.696 -> -2396.5586
$ws.Cells.Item(110, 14).Value = -6325.1428  # ARM N110: -6752.4614 -> -6325.1428

$ws.Cells.Item(123, 8).Value = 33125  # ARM H123: 54385.8 -> 33125
$ws.Cells.Item(123, 10).Value = 33125  # ARM J123: 54385.8 -> 33125
$ws.Cells.Item(123, 12).Value = 33125  # ARM L123: 54385.8 -> 33125
$ws.Cells.Item(123, 14).Value = -42925  # ARM N123: -64185.8 -> -42925

$ws.Cells.Item(132, 8).Value = 1512.6482  # ARM H132: 1601.66 -> 1512.6482
$ws.Cells.Item(132, 9).Value = 1197  # ARM I132: 1293.8857 -> 1197
$ws.Cells.Item(132, 10).Value = 2262.3125  # ARM J132: 2319.8 -> 2262.3125
$ws.Cells.Item(132, 11).Value = 3591  # ARM K132: 3881.6571 -> 3591
$ws.Cells.Item(132, 12).Value = 6786.9375  # ARM L132: 6959.400000000001 -> 6786.9375
$ws.Cells.Item(132, 13).Value = -1061  # ARM M132: -1351.6571 -> -1061
$ws.Cells.Item(132, 14).Value = -11846.9375  # ARM N132: -12019.4 -> -11846.9375

$ws.Cells.Item(136, 8).Value = 1152.303  # ARM H136: 1208.4667 -> 1152.303
$ws.Cells.Item(136, 9).Value = 1003.9048  # ARM I136: 1098.8235 -> 1003.9048
$ws.Cells.Item(136, 10).Value = 1412  # ARM J136: 1351.8462 -> 1412
$ws.Cells.Item(136, 11).Value = 3011.7144  # ARM K136: 3296.4705 -> 3011.7144
$ws.Cells.Item(136, 12).Value = 4236  # ARM L136: 4055.5386 -> 4236
$ws.Cells.Item(136, 13).Value = -461.7143999999998  # ARM M136: -746.4704999999999 -> -461.7143999999998
$ws.Cells.Item(136, 14).Value = -9336  # ARM N136: -9155.5386 -> -9336

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(2, 8).Value = 22680  # BSM H2: 24925 -> 22680
$ws.Cells.Item(2, 10).Value = 22680  # BSM J2: 24925 -> 22680
$ws.Cells.Item(2, 12).Value = 22680  # BSM L2: 24925 -> 22680
$ws.Cells.Item(2, 14).Value = -22906  # BSM N2: -25151 -> -22906

$ws.Cells.Item(86, 8).Value = 16676034  # BSM H86: 20011020 -> 16676034
$ws.Cells.Item(86, 9).Value = 22223468  # BSM I86: 28572716 -> 22223468
$ws.Cells.Item(86, 11).Value = 22223468  # BSM K86: 28572716 -> 22223468
$ws.Cells.Item(86, 13).Value = -22222345  # BSM M86: -28571593 -> -22222345

$ws.Cells.Item(89, 8).Value = 16676034  # BSM H89: 20011020 -> 16676034
$ws.Cells.Item(89, 9).Value = 22223468  # BSM I89: 28572716 -> 22223468
$ws.Cells.Item(89, 11).Value = 111117340  # BSM K89: 142863580 -> 111117340
$ws.Cells.Item(89, 13).Value = -111111724  # BSM M89: -142857964 -> -111111724

$ws.Cells.Item(94, 8).Value = 18003.25  # BSM H94: 10417.477 -> 18003.25
$ws.Cells.Item(94, 9).Value = 1166.1666  # BSM I94: 775.3125 -> 1166.1666
$ws.Cells.Item(94, 10).Value = 34840.332  # BSM J94: 41272.4 -> 34840.332
$ws.Cells.Item(94, 11).Value = 1166.1666  # BSM K94: 775.3125 -> 1166.1666
$ws.Cells.Item(94, 12).Value = 34840.332  # BSM L94: 41272.4 -> 34840.332
$ws.Cells.Item(94, 13).Value = -715.1666  # BSM M94: -324.3125 -> -715.1666
$ws.Cells.Item(94, 14).Value = -35742.332  # BSM N94: -42174.4 -> -35742.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 29414224  # CRP H16: 41669308 -> 29414224
$ws.Cells.Item(16, 9).Value = 38464024  # CRP I16: 45457020 -> 38464024
$ws.Cells.Item(16, 10).Value = 2372.5  # CRP J16: 4500 -> 2372.5
$ws.Cells.Item(16, 11).Value = 38464024  # CRP K16: 45457020 -> 38464024
$ws.Cells.Item(16, 12).Value = 2372.5  # CRP L16: 4500 -> 2372.5
$ws.Cells.Item(16, 13).Value = -38463737  # CRP M16: -45456733 -> -38463737
$ws.Cells.Item(16, 14).Value = -2946.5  # CRP N16: -5074 -> -2946.5

$ws.Cells.Item(86, 8).Value = 253525.7  # CRP H86: 316374.5 -> 253525.7
$ws.Cells.Item(86, 9).Value = 360907.56  # CRP I86: 420516.16 -> 360907.56
$ws.Cells.Item(86, 10).Value = 2968  # CRP J86: 3949.5 -> 2968
$ws.Cells.Item(86, 11).Value = 360907.56  # CRP K86: 420516.16 -> 360907.56
$ws.Cells.Item(86, 12).Value = 2968  # CRP L86: 3949.5 -> 2968
$ws.Cells.Item(86, 13).Value = -359784.56  # CRP M86: -419393.16 -> -359784.56
$ws.Cells.Item(86, 14).Value = -5214  # CRP N86: -6195.5 -> -5214

$ws.Cells.Item(89, 8).Value = 253525.7  # CRP H89: 316374.5 -> 253525.7
$ws.Cells.Item(89, 9).Value = 360907.56  # CRP I89: 420516.16 -> 360907.56
$ws.Cells.Item(89, 10).Value = 2968  # CRP J89: 3949.5 -> 2968
$ws.Cells.Item(89, 11).Value = 1804537.8  # CRP K89: 2102580.8 -> 1804537.8
$ws.Cells.Item(89, 12).Value = 14840  # CRP L89: 19747.5 -> 14840
$ws.Cells.Item(89, 13).Value = -1798921.8  # CRP M89: -2096964.8 -> -1798921.8
$ws.Cells.Item(89, 14).Value = -26072  # CRP N89: -30979.5 -> -26072

$ws.Cells.Item(105, 8).Value = 698.75  # CRP H105: 667.4167 -> 698.75
$ws.Cells.Item(105, 9).Value = 590  # CRP I105: 520 -> 590
$ws.Cells.Item(105, 10).Value = 880  # CRP J105: 873.8 -> 880
$ws.Cells.Item(105, 11).Value = 590  # CRP K105: 520 -> 590
$ws.Cells.Item(105, 12).Value = 880  # CRP L105: 873.8 -> 880
$ws.Cells.Item(105, 13).Value = 1157  # CRP M105: 1227 -> 1157
$ws.Cells.Item(105, 14).Value = -4374  # CRP N105: -4367.8 -> -4374

$ws.Cells.Item(107, 8).Value = 25001226  # CRP H107: 22728418 -> 25001226
$ws.Cells.Item(107, 9).Value = 41667708  # CRP I107: 35715230 -> 41667708
$ws.Cells.Item(107, 10).Value = 1503.375  # CRP J107: 1498.875 -> 1503.375
$ws.Cells.Item(107, 11).Value = 41667708  # CRP K107: 35715230 -> 41667708
$ws.Cells.Item(107, 12).Value = 1503.375  # CRP L107: 1498.875 -> 1503.375
$ws.Cells.Item(107, 13).Value = -41665788  # CRP M107: -35713310 -> -41665788
$ws.Cells.Item(107, 14).Value = -5343.375  # CRP N107: -5338.875 -> -5343.375

$ws.Cells.Item(113, 8).Value = 29414224  # CRP H113: 41669308 -> 29414224
$ws.Cells.Item(113, 9).Value = 38464024  # CRP I113: 45457020 -> 38464024
$ws.Cells.Item(113, 10).Value = 2372.5  # CRP J113: 4500 -> 2372.5
$ws.Cells.Item(113, 11).Value = 38464024  # CRP K113: 45457020 -> 38464024
$ws.Cells.Item(113, 12).Value = 2372.5  # CRP L113: 4500 -> 2372.5
$ws.Cells.Item(113, 13).Value = -38461854  # CRP M113: -45454850 -> -38461854
$ws.Cells.Item(113, 14).Value = -6712.5  # CRP N113: -8840 -> -6712.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 776.875  # CUL H7: 957.5333000000001 -> 776.875
$ws.Cells.Item(7, 9).Value = 147.77777  # CUL I7: 267.44446 -> 147.77777
$ws.Cells.Item(7, 10).Value = 1585.7142  # CUL J7: 1992.6666 -> 1585.7142
$ws.Cells.Item(7, 11).Value = 443.33331  # CUL K7: 802.33338 -> 443.33331
$ws.Cells.Item(7, 12).Value = 4757.142599999999  # CUL L7: 5977.9998 -> 4757.142599999999
$ws.Cells.Item(7, 13).Value = -331.33331  # CUL M7: -690.33338 -> -331.33331
$ws.Cells.Item(7, 14).Value = -4981.142599999999  # CUL N7: -6201.9998 -> -4981.142599999999

$ws.Cells.Item(131, 8).Value = 930.9388  # CUL H131: 944.3778 -> 930.9388
$ws.Cells.Item(131, 9).Value = 358  # CUL I131: 370.92307 -> 358
$ws.Cells.Item(131, 10).Value = 1160.1143  # CUL J131: 1177.3438 -> 1160.1143
$ws.Cells.Item(131, 11).Value = 1074  # CUL K131: 1112.76921 -> 1074
$ws.Cells.Item(131, 12).Value = 3480.3429  # CUL L131: 3532.0314 -> 3480.3429
$ws.Cells.Item(131, 13).Value = 3966  # CUL M131: 3927.23079 -> 3966
$ws.Cells.Item(131, 14).Value = -13560.3429  # CUL N131: -13612.0314 -> -13560.3429

$ws.Cells.Item(140, 8).Value = 1414.138  # CUL H140: 1774.4445 -> 1414.138
$ws.Cells.Item(140, 9).Value = 1630.909  # CUL I140: 1664 -> 1630.909
$ws.Cells.Item(140, 10).Value = 1281.6666  # CUL J140: 1912.5 -> 1281.6666
$ws.Cells.Item(140, 11).Value = 4892.727000000001  # CUL K140: 4992 -> 4892.727000000001
$ws.Cells.Item(140, 12).Value = 3844.9998  # CUL L140: 5737.5 -> 3844.9998
$ws.Cells.Item(140, 13).Value = 287.2729999999992  # CUL M140: 188 -> 287.2729999999992
$ws.Cells.Item(140, 14).Value = -14204.9998  # CUL N140: -16097.5 -> -14204.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(19, 8).Value = 0  # GSM H19: 405 -> 0
$ws.Cells.Item(19, 9).Value = 0  # GSM I19: 405 -> 0
$ws.Cells.Item(19, 11).Value = 0  # GSM K19: 405 -> 0
$ws.Cells.Item(19, 13).ClearContents()  # GSM M19: remove (was -117)

$ws.Cells.Item(102, 8).Value = 1609.28  # GSM H102: 1744.4231 -> 1609.28
$ws.Cells.Item(102, 9).Value = 1781.7  # GSM I102: 1814.7916 -> 1781.7
$ws.Cells.Item(102, 10).Value = 919.6  # GSM J102: 900 -> 919.6
$ws.Cells.Item(102, 11).Value = 1781.7  # GSM K102: 1814.7916 -> 1781.7
$ws.Cells.Item(102, 12).Value = 919.6  # GSM L102: 900 -> 919.6
$ws.Cells.Item(102, 13).Value = -159.7  # GSM M102: -192.7916 -> -159.7
$ws.Cells.Item(102, 14).Value = -4163.6  # GSM N102: -4144 -> -4163.6

$ws.Cells.Item(123, 8).Value = 17078.75  # GSM H123: 10025.4 -> 17078.75
$ws.Cells.Item(123, 10).Value = 17078.75  # GSM J123: 10025.4 -> 17078.75
$ws.Cells.Item(123, 12).Value = 17078.75  # GSM L123: 10025.4 -> 17078.75
$ws.Cells.Item(123, 14).Value = -21978.75  # GSM N123: -14925.4 -> -21978.75

$ws.Cells.Item(126, 8).Value = 2919.45  # GSM H126: 3140.5881 -> 2919.45
$ws.Cells.Item(126, 9).Value = 1954.3334  # GSM I126: 1998.5714 -> 1954.3334
$ws.Cells.Item(126, 10).Value = 3709.0908  # GSM J126: 3940 -> 3709.0908
$ws.Cells.Item(126, 11).Value = 5863.0002  # GSM K126: 5995.7142 -> 5863.0002
$ws.Cells.Item(126, 12).Value = 11127.2724  # GSM L126: 11820 -> 11127.2724
$ws.Cells.Item(126, 13).Value = -3393.0002  # GSM M126: -3525.7142 -> -3393.0002
$ws.Cells.Item(126, 14).Value = -16067.2724  # GSM N126: -16760 -> -16067.2724

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 359.33334  # LTW H55: 306.05884 -> 359.33334
$ws.Cells.Item(55, 9).Value = 384.2857  # LTW I55: 274.3 -> 384.2857
$ws.Cells.Item(55, 10).Value = 337.5  # LTW J55: 351.42856 -> 337.5
$ws.Cells.Item(55, 11).Value = 384.2857  # LTW K55: 274.3 -> 384.2857
$ws.Cells.Item(55, 12).Value = 337.5  # LTW L55: 351.42856 -> 337.5
$ws.Cells.Item(55, 13).Value = -211.2857  # LTW M55: -101.3 -> -211.2857
$ws.Cells.Item(55, 14).Value = -683.5  # LTW N55: -697.4285600000001 -> -683.5

$ws.Cells.Item(117, 8).Value = 27496  # LTW H117: 36000 -> 27496
$ws.Cells.Item(117, 10).Value = 27496  # LTW J117: 36000 -> 27496
$ws.Cells.Item(117, 12).Value = 27496  # LTW L117: 36000 -> 27496
$ws.Cells.Item(117, 14).Value = -36674  # LTW N117: -45178 -> -36674

$ws.Cells.Item(122, 8).Value = 5680  # LTW H122: 5966.6665 -> 5680
$ws.Cells.Item(122, 9).Value = 5722.222  # LTW I122: 6100 -> 5722.222
$ws.Cells.Item(122, 11).Value = 17166.666  # LTW K122: 18300 -> 17166.666
$ws.Cells.Item(122, 13).Value = -14716.666  # LTW M122: -15850 -> -14716.666

$ws.Cells.Item(132, 8).Value = 6871765.5  # LTW H132: 7026019.5 -> 6871765.5
$ws.Cells.Item(132, 9).Value = 13026051  # LTW I132: 13894159 -> 13026051
$ws.Cells.Item(132, 10).Value = 1865.0698  # LTW J132: 1786.7273 -> 1865.0698
$ws.Cells.Item(132, 11).Value = 39078153  # LTW K132: 41682477 -> 39078153
$ws.Cells.Item(132, 12).Value = 5595.2094  # LTW L132: 5360.1819 -> 5595.2094
$ws.Cells.Item(132, 13).Value = -39075623  # LTW M132: -41679947 -> -39075623
$ws.Cells.Item(132, 14).Value = -10655.2094  # LTW N132: -10420.1819 -> -10655.2094

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 100  # WVR H14: 0 -> 100
$ws.Cells.Item(14, 10).Value = 100  # WVR J14: 0 -> 100
$ws.Cells.Item(14, 12).Value = 100  # WVR L14: 0 -> 100
$ws.Cells.Item(14, 14).Value = -436  # WVR N14: None -> -436

$ws.Cells.Item(118, 8).Value = 39695  # WVR H118: 40392 -> 39695
$ws.Cells.Item(118, 10).Value = 39695  # WVR J118: 40392 -> 39695
$ws.Cells.Item(118, 12).Value = 39695  # WVR L118: 40392 -> 39695
$ws.Cells.Item(118, 14).Value = -43009  # WVR N118: -43706 -> -43009

$ws.Cells.Item(121, 8).Value = 29000  # WVR H121: 0 -> 29000
$ws.Cells.Item(121, 10).Value = 29000  # WVR J121: 0 -> 29000
$ws.Cells.Item(121, 12).Value = 29000  # WVR L121: 0 -> 29000
$ws.Cells.Item(121, 14).Value = -32494  # WVR N121: None -> -32494

$ws.Cells.Item(122, 8).Value = 1650.3334  # WVR H122: 1430.5238 -> 1650.3334
$ws.Cells.Item(122, 9).Value = 1264.3636  # WVR I122: 1127.625 -> 1264.3636
$ws.Cells.Item(122, 10).Value = 2256.8572  # WVR J122: 2399.8 -> 2256.8572
$ws.Cells.Item(122, 11).Value = 3793.0908  # WVR K122: 3382.875 -> 3793.0908
$ws.Cells.Item(122, 12).Value = 6770.571599999999  # WVR L122: 7199.400000000001 -> 6770.571599999999
$ws.Cells.Item(122, 13).Value = -1343.0908  # WVR M122: -932.875 -> -1343.0908
$ws.Cells.Item(122, 14).Value = -11670.5716  # WVR N122: -12099.4 -> -11670.5716

$ws.Cells.Item(123, 8).Value = 38590.727  # WVR H123: 44315.207 -> 38590.727
$ws.Cells.Item(123, 10).Value = 38590.727  # WVR J123: 44315.207 -> 38590.727
$ws.Cells.Item(123, 12).Value = 38590.727  # WVR L123: 44315.207 -> 38590.727
$ws.Cells.Item(123, 14).Value = -48390.727  # WVR N123: -54115.207 -> -48390.727
